$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B for "status_label"
$ws.Columns.Item(2).Insert()

$arr = New-Object 'object[,]' 20,10

$arr[0,0] = "statut"
$arr[0,1] = "status_label"
$arr[0,2] = "NCTId"
$arr[0,3] = "eudraCT"
$arr[0,4] = "completion_year"
$arr[0,5] = "clinical_trial_title"
$arr[0,6] = "acronym"
$arr[0,7] = "results_1y"
$arr[0,8] = "results_3y"
$arr[0,9] = "results"
$arr[1,0] = "🟥"
$arr[1,1] = "rouge"
$arr[1,2] = "NCT01448018"
$arr[1,3] = ""
$arr[1,4] = "2013"
$arr[1,5] = "Pilot Study on Efficacy and Tolerance of Intravitreous Injection of Ranibizumab (Lucentis®) in Early-onset Central Retinal Vein Occlusion in Comparison to Hemodilution Using Erythrocytapheresis"
$arr[1,6] = "CHIC-3"
$arr[1,7] = $false
$arr[1,8] = $false
$arr[1,9] = $false
$arr[2,0] = "🟥"
$arr[2,1] = "rouge"
$arr[2,2] = "NCT02157077"
$arr[2,3] = ""
$arr[2,4] = "2015"
$arr[2,5] = "A Phase III b, Multicenter Study of the Efficacy and Safety of Aflibercept Switch in Patients With Exudative AMD With Detachment of the Retinal Pigment Epithelium and Previously Treated With Ranibizumab Intravitreal Injection. (ARI2)"
$arr[2,6] = "ARI2"
$arr[2,7] = $false
$arr[2,8] = $false
$arr[2,9] = $true
$arr[3,0] = "🟧"
$arr[3,1] = "orange"
$arr[3,2] = "NCT02700893"
$arr[3,3] = ""
$arr[3,4] = "2016"
$arr[3,5] = "Cerebral NIRS Profiles During Premedication for Neonatal Intubation"
$arr[3,6] = ""
$arr[3,7] = $false
$arr[3,8] = $true
$arr[3,9] = $true
$arr[4,0] = "🟥"
$arr[4,1] = "rouge"
$arr[4,2] = "NCT02806830"
$arr[4,3] = ""
$arr[4,4] = "2017"
$arr[4,5] = "Evaluation de la gêne Oculaire après Injections intravitréennes"
$arr[4,6] = "EVAGO"
$arr[4,7] = $false
$arr[4,8] = $false
$arr[4,9] = $false
$arr[5,0] = "🟥"
$arr[5,1] = "rouge"
$arr[5,2] = "NCT02899806"
$arr[5,3] = ""
$arr[5,4] = "2017"
$arr[5,5] = "Impact of a Video Explaining Epidural Analgesia in Obstetrics in Terms of Satisfaction, Understanding and Anxiety: A Prospective Randomised Trial"
$arr[5,6] = "VIDEOCLIP"
$arr[5,7] = $false
$arr[5,8] = $false
$arr[5,9] = $false
$arr[6,0] = "🟩"
$arr[6,1] = "vert"
$arr[6,2] = "NCT04318431"
$arr[6,3] = ""
$arr[6,4] = "2020"
$arr[6,5] = "Prevalence of SARS -Cov2 Carriage in Asymptomatic and Mildly-symptomatic Children, a Cross-sectional, Prospective, Multicentre, Observational Study in Primary Care."
$arr[6,6] = "COVILLE"
$arr[6,7] = $true
$arr[6,8] = $true
$arr[6,9] = $true
$arr[7,0] = "🟥"
$arr[7,1] = "rouge"
$arr[7,2] = "NCT04583189"
$arr[7,3] = ""
$arr[7,4] = "2020"
$arr[7,5] = "Evaluation Des Performances du Test Rapide antigénique Covid-19 Ag BSS Chez l'Enfant Symptomatique Dans un Service d'Urgences pédiatriques"
$arr[7,6] = ""
$arr[7,7] = $false
$arr[7,8] = $false
$arr[7,9] = $false
$arr[8,0] = "🟩"
$arr[8,1] = "vert"
$arr[8,2] = "NCT01490580"
$arr[8,3] = ""
$arr[8,4] = "2020"
$arr[8,5] = "Double Blind Randomized Controlled Trial Comparing `"Atropine+Propofol`" Versus `"Atropine+Atracurium+Sufentanil`" as a Premedication Prior to Semi-urgent or Elective Endotracheal Intubation of Term and Preterm Newborns"
$arr[8,6] = "PRETTINEO"
$arr[8,7] = $true
$arr[8,8] = $true
$arr[8,9] = $true
$arr[9,0] = "🟥"
$arr[9,1] = "rouge"
$arr[9,2] = "NCT03488758"
$arr[9,3] = ""
$arr[9,4] = "2021"
$arr[9,5] = "Enjoyment of Infant Formulas Based on Cow or Goat Milk Protein"
$arr[9,6] = "CHARLIE"
$arr[9,7] = $false
$arr[9,8] = $false
$arr[9,9] = $false
$arr[10,0] = "🟥"
$arr[10,1] = "rouge"
$arr[10,2] = "NCT04776174"
$arr[10,3] = ""
$arr[10,4] = "2021"
$arr[10,5] = "Crossover Comparison of the Efficacy and Tolerance of Telerobotic vs Standard Ultrasound Exam in Children"
$arr[10,6] = ""
$arr[10,7] = $false
$arr[10,8] = $false
$arr[10,9] = $false
$arr[11,0] = "🟥"
$arr[11,1] = "rouge"
$arr[11,2] = "NCT03939377"
$arr[11,3] = ""
$arr[11,4] = "2021"
$arr[11,5] = "Evaluation of an Osteopathic Procedure in the Management of Pain in Palliative Care Patients in a Mobile Palliative Care Support Team (EMASP): Controlled, Randomized, Single-blind Study"
$arr[11,6] = "OSTEOPAL"
$arr[11,7] = $false
$arr[11,8] = $false
$arr[11,9] = $false
$arr[12,0] = "🟧"
$arr[12,1] = "orange"
$arr[12,2] = "NCT03030664"
$arr[12,3] = ""
$arr[12,4] = "2021"
$arr[12,5] = "Randomised Controlled Trial With Two Parallel Arms Testing the Effect of L. Reuteri on Bowel Movements in Children Aged 6 Months to 4 Years"
$arr[12,6] = "BIOWELL"
$arr[12,7] = $false
$arr[12,8] = $true
$arr[12,9] = $true
$arr[13,0] = "🟥"
$arr[13,1] = "rouge"
$arr[13,2] = "NCT03803228"
$arr[13,3] = ""
$arr[13,4] = "2021"
$arr[13,5] = "Comparison of the Cumulative Number of Oocytes Obtained With 2 Controlled Ovarian Hyperstimulations (COH) Within the Same Cycle With FertistartKit® (DUOSTIM) Versus 2 Conventional COH in Poor Ovarian Responders Undergoing IVF. Bistim Study"
$arr[13,6] = ""
$arr[13,7] = $false
$arr[13,8] = $false
$arr[13,9] = $false
$arr[14,0] = "🟥"
$arr[14,1] = "rouge"
$arr[14,2] = "NCT04068558"
$arr[14,3] = ""
$arr[14,4] = "2021"
$arr[14,5] = "Synchronized Nasal Intermittent Positive Pressure Ventilation Versus Noninvasive Neurally Adjusted Ventilatory Assist Ventilation in Extremely Premature Infants: a Randomized Crossover Trial"
$arr[14,6] = "EASYNNEO"
$arr[14,7] = $false
$arr[14,8] = $false
$arr[14,9] = $false
$arr[15,0] = "🟥"
$arr[15,1] = "rouge"
$arr[15,2] = "NCT05079139"
$arr[15,3] = ""
$arr[15,4] = "2022"
$arr[15,5] = "Musset's Surgical Technique: Evaluation of Long-term Results (LONGOMUSSET)"
$arr[15,6] = "LONGOMUSSET"
$arr[15,7] = $false
$arr[15,8] = $false
$arr[15,9] = $false
$arr[16,0] = "🟥"
$arr[16,1] = "rouge"
$arr[16,2] = "NCT03895099"
$arr[16,3] = ""
$arr[16,4] = "2023"
$arr[16,5] = "Feasibility and Efficacy of a New Ovarian Stimulation Regimen With RANDom Start, Use of Corifollitropin Alpha and Progestin Protocol for Oocyte donorS"
$arr[16,6] = "RANDOS"
$arr[16,7] = $false
$arr[16,8] = $false
$arr[16,9] = $false
$arr[17,0] = "🟥"
$arr[17,1] = "rouge"
$arr[17,2] = "NCT03540706"
$arr[17,3] = ""
$arr[17,4] = "2023"
$arr[17,5] = "Impact of the Use of C-reactive Protein in a Micro-method on the Prescription of Antibiotics in General Practitioners Consulting in the Office"
$arr[17,6] = "VIP"
$arr[17,7] = $false
$arr[17,8] = $false
$arr[17,9] = $false
$arr[18,0] = "🟥"
$arr[18,1] = "rouge"
$arr[18,2] = "NCT02884245"
$arr[18,3] = ""
$arr[18,4] = "2023"
$arr[18,5] = "Interest of Estrogen Scheduling Before Ovarian Stimulation With Corifollitropin Alfa in Women Older Than 38 Years Old Undergoing in Vitro Fertilization"
$arr[18,6] = "PRESCORI"
$arr[18,7] = $false
$arr[18,8] = $false
$arr[18,9] = $false
$arr[19,0] = "🟥"
$arr[19,1] = "rouge"
$arr[19,2] = "NCT04667065"
$arr[19,3] = ""
$arr[19,4] = "2023"
$arr[19,5] = "Evaluation of a Remotely Guided Physical Preparation by a Physical Activity Teacher Adapted With the Help of a Smartwatch Before Bronchial Cancer Surgery"
$arr[19,6] = "PREPACHIR"
$arr[19,7] = $false
$arr[19,8] = $false
$arr[19,9] = $false

$ws.Range("A1:J20").Value = $arr
